# BC Hydro Survey Helper - final draft edits
#
# Strategy: locate each paragraph that needs editing via Find (unique
# substring), expand to the whole paragraph, and replace its contents with
# InsertXML using a minimal WordprocessingML fragment. This lets us merge
# runs and drop/relocate proofErr + bookmark markers precisely, matching
# the target OOXML, without depending on brittle paragraph indices.

$d = $word.ActiveDocument

function Replace-ParagraphByFind($searchText, $xmlFragment) {
    $rng = $d.Content.Duplicate
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $searchText"
    }
    $para = $rng.Paragraphs(1).Range
    $wrapped = '<?xml version="1.0" standalone="yes"?>' + `
        '<?mso-application progid="Word.Document"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body>' + $xmlFragment + '</w:body></w:document>' + `
        '</pkg:xmlData></pkg:part></pkg:package>'
    [void]$para.InsertXML($wrapped)
}

# 1) "For Dr. Jahangir Hossain, Ph.D" -> merge runs, drop proofErr, add "." run
Replace-ParagraphByFind "For Dr. Jahangir Hossain, Ph.D" (
    '<w:p><w:pPr><w:pStyle w:val="Subtitle"/><w:jc w:val="center"/></w:pPr>' +
    '<w:r><w:t>For Dr. Jahangir Hossain, Ph.D</w:t></w:r>' +
    '<w:r><w:t>.</w:t></w:r></w:p>'
)

# 2) "Amec Foster Wheeler plc" intro paragraph -> merge runs, drop proofErr
Replace-ParagraphByFind "Amec Foster Wheeler" (
    '<w:p><w:r><w:tab/><w:t xml:space="preserve">This document is a discussion of the structure and design of the as-of-yet-unnamed application to be developed for </w:t></w:r>' +
    '<w:r><w:t>Amec Foster Wheeler plc</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> (Amec from here on) to assist in the process of surveying BC Hydro power transmission poles. The goal of this discussion is to relate an early concept for the overall </w:t></w:r>' +
    '<w:r><w:t>structure</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> of the application, as well as to reveal some of the underlying </w:t></w:r>' +
    '<w:r><w:t>framework needed to support it.</w:t></w:r></w:p>'
)

# 3) "Workflow Overview" body paragraph -> merge runs, drop proofErr (gramStart/End)
Replace-ParagraphByFind "The application can be considered as two flows" (
    '<w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">The application can be considered as two flows. The UI Flow experienced by the user, and the flow of data handled behind the scenes. In the diagram below, each block shows a piece of the app that must be developed. This preliminary layout shows how the user will login, view a list of assigned poles to survey, then complete a the set of tasks associated with that assignment. Once the survey is complete it can b reviewed before either placing in a local data repository or, if network access is available, sending it directly to remote storage and processing. Each of the elements of this workflow are described in brief in the next section.</w:t></w:r></w:p>'
)

# 4) "Workflow Components" heading -> drop the _GoBack bookmark pair
Replace-ParagraphByFind "Workflow Components" (
    '<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr>' +
    '<w:r><w:lastRenderedPageBreak/><w:t>Workflow Components</w:t></w:r></w:p>'
)

# 5) "Poles may have a set of documents" paragraph -> merge runs, drop proofErr
Replace-ParagraphByFind "Poles may have a set of documents" (
    '<w:p><w:r><w:tab/></w:r>' +
    '<w:r><w:t>Poles may have a set of documents associated with them, including maps, safety, environmental and legal information. This page provides the ability to select and view these documents one by one, and requests confirmation from the user that all documents have been reviewed before proceeding. These documents are provided by AMEC as part of the aforementioned pole assignments.</w:t></w:r></w:p>'
)

# 6) "Data collected by the user" paragraph -> merge runs, drop proofErr, and
#    relocate the _GoBack bookmark pair into the middle of "it's" ("it" | "s")
Replace-ParagraphByFind "Data collected by the user" (
    '<w:p><w:r><w:tab/></w:r>' +
    '<w:r><w:t>Data collected by the user is either stored in a local repository using the device' + [char]0x2019 + 's filesystem, or if network access is available sent immediately to remote servers for storage and processing. Data may be checksummed and/or encrypted locally to ensure it</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '<w:r><w:t>s integrity and security.</w:t></w:r></w:p>'
)

Write-Output "All replacements applied."
